# Wrap the QIDCELL when Error Handling
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- N4: update the QID text value (keep it stored as text, not a number) ---
$n4 = $ws.Cells.Item(4, 14)
$n4.Formula = '="412051"'
$n4.Copy()
$n4.PasteSpecial(-4163)  # xlPasteValues

# --- N6: update the QID text value (keep it stored as text, not a number) ---
$n6 = $ws.Cells.Item(6, 14)
$n6.Formula = '="412050"'
$n6.Copy()
$n6.PasteSpecial(-4163)  # xlPasteValues

# --- N5: update error message text, and wrap + top-align it (new cell style) ---
$n5 = $ws.Cells.Item(5, 14)
$n5.Value = "Question Failed To Create:Object failed to Match any listed under this product"
$n5.WrapText = $true
$n5.VerticalAlignment = -4160  # xlTop

$excel.CutCopyMode = $false
